$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterData")

# --- Header row: bold the existing header cells A1:L1 ---
$ws.Range("A1:L1").Font.Bold = $true

# --- New "Account No" header cell M1 (bold, matches existing bold font used elsewhere) ---
$ws.Range("M1").Value = "Account No"
$ws.Range("M1").Font.Bold = $true

# --- New row 6 data: currentUserRegisterToLogin test case ---
$ws.Range("H6").Value = 5111111111
$ws.Range("I6").Value = "testautomation7295@gmail.com"
$ws.Range("J6").Value = "Automation1212"
$ws.Range("K6").Value = "Register sonrası login başarısız"
$ws.Range("M6").Value = 22501947

# Give I6 the same visual style (font) as the other email/hyperlink cells before linking it
$ws.Range("I2").Copy()
$ws.Range("I6").PasteSpecial(-4122)

# Hyperlink for the new email cell
$ws.Hyperlinks.Add($ws.Range("I6"), "mailto:testautomation7295@gmail.com", "", "", "testautomation7295@gmail.com")

# Restore the matching (non-hyperlink-builtin) style that Hyperlinks.Add just overwrote
$ws.Range("I2").Copy()
$ws.Range("I6").PasteSpecial(-4122)

# --- View: selection / scroll position moved to reflect the newly visible columns ---
$ws.Range("H10").Select()
